# Updates the crypto price/volume table (Sheet1) produced by the
# "Updated cryptos list" GitHub Actions job.
#
# Column D ("Price") holds values that are stored as literal text in the
# workbook (e.g. "152.40", "1.000", "27.344.87") rather than numbers, so
# any new value that LOOKS like a plain decimal number (and would
# otherwise be auto-converted to a true numeric cell by Excel) is written
# with the cell pre-formatted as Text ("@") to keep it a literal string -
# this preserves things like trailing zeros ("103.00") that a numeric
# cell cannot represent. Column E ("Volume(1h)") values already contain
# a "%" sign and padding spaces, so Excel never mistakes them for numbers
# and they can be written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to Text format first so a numeric-looking string
    # (e.g. "1.002", "0.000009066", "103.00") is stored verbatim instead
    # of being parsed into a Double by Excel.
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
}


# Row 2
$ws.Range("D2").Value = '27.344.87'

# Row 3
$ws.Range("D3").Value = '1.857.08'
$ws.Range("E3").Value = '  +1.58%  '

# Row 4
Set-TextValue "D4" '1.002'
$ws.Range("E4").Value = '  -0.71%  '

# Row 5
Set-TextValue "D5" '314.24'
$ws.Range("E5").Value = '  +0.93%  '

# Row 6
Set-TextValue "D6" '1.001'
$ws.Range("E6").Value = '  -0.67%  '

# Row 7
Set-TextValue "D7" '0.4614'
$ws.Range("E7").Value = '  -0.85%  '

# Row 8
$ws.Range("E8").Value = '  +0.35%  '

# Row 9
Set-TextValue "D9" '0.07323'
$ws.Range("E9").Value = '  -0.46%  '

# Row 10
Set-TextValue "D10" '0.8812'
$ws.Range("E10").Value = '  +0.91%  '

# Row 11
$ws.Range("E11").Value = '  -0.09%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue "D12" '0.07804'
$ws.Range("E12").Value = '  -0.89%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.922.16'
$ws.Range("E13").Value = '  +6.51%  '

# Row 14
$ws.Range("E14").Value = '  +0.61%  '

# Row 15
Set-TextValue "D15" '6.552'
$ws.Range("E15").Value = '  -0.31%  '

# Row 16
Set-TextValue "D16" '91.88'
$ws.Range("E16").Value = '  -0.07%  '

# Row 17
Set-TextValue "D17" '1.002'
$ws.Range("E17").Value = '  -0.72%  '

# Row 18
Set-TextValue "D18" '0.000009066'
$ws.Range("E18").Value = '  +2.23%  '

# Row 19
$ws.Range("E19").Value = '  -0.63%  '

# Row 20
Set-TextValue "D20" '14.79'
$ws.Range("E20").Value = '  +0.64%  '

# Row 21
$ws.Range("D21").Value = '27.355.82'
$ws.Range("E21").Value = '  +1.95%  '

# Row 22
Set-TextValue "D22" '5.130'
$ws.Range("E22").Value = '  -0.46%  '

# Row 23
Set-TextValue "D23" '10.52'
$ws.Range("E23").Value = '  -0.40%  '

# Row 24
$ws.Range("D24").Value = '2.153.45'
$ws.Range("E24").Value = '  +2.96%  '

# Row 25
$ws.Range("E25").Value = '  +5.39%  '

# Row 26
Set-TextValue "D26" '152.40'
$ws.Range("E26").Value = '  -0.12%  '

# Row 27
$ws.Range("E27").Value = '  +0.91%  '

# Row 28
Set-TextValue "D28" '2.073'
$ws.Range("E28").Value = '  -1.08%  '

# Row 29
Set-TextValue "D29" '5.106'
$ws.Range("E29").Value = '  -0.35%  '

# Row 30
Set-TextValue "D30" '116.09'
$ws.Range("E30").Value = '  +0.50%  '

# Row 31
Set-TextValue "D31" '0.08863'
$ws.Range("E31").Value = '  -0.09%  '

# Row 32
Set-TextValue "D32" '0.7710'
$ws.Range("E32").Value = '  +6.00%  '

# Row 33
Set-TextValue "D33" '3.041'
$ws.Range("E33").Value = '  +1.88%  '

# Row 34
$ws.Range("E34").Value = '  +3.30%  '

# Row 35
Set-TextValue "D35" '4.495'
$ws.Range("E35").Value = '  +1.22%  '

# Row 36
Set-TextValue "D36" '2.648'
$ws.Range("E36").Value = '  +5.98%  '

# Row 37
$ws.Range("E37").Value = '  +0.05%  '

# Row 38
$ws.Range("E38").Value = '  +0.52%  '

# Row 39
Set-TextValue "D39" '0.05234'
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
Set-TextValue "D40" '2.954'
$ws.Range("E40").Value = '  +0.84%  '

# Row 41
Set-TextValue "D41" '7.019'
$ws.Range("E41").Value = '  -3.58%  '

# Row 42
$ws.Range("E42").Value = '  -0.86%  '

# Row 43
Set-TextValue "D43" '0.1637'
$ws.Range("E43").Value = '  +0.81%  '

# Row 44
Set-TextValue "D44" '8.409'
$ws.Range("E44").Value = '  +2.42%  '

# Row 45
Set-TextValue "D45" '0.4823'
$ws.Range("E45").Value = '  -0.29%  '

# Row 46
Set-TextValue "D46" '10.32'
$ws.Range("E46").Value = '  +1.44%  '

# Row 47
Set-TextValue "D47" '1.001'
$ws.Range("E47").Value = '  -0.75%  '

# Row 48
Set-TextValue "D48" '103.00'
$ws.Range("E48").Value = '  +0.41%  '

# Row 49
$ws.Range("E49").Value = '  +1.80%  '

# Row 50
Set-TextValue "D50" '0.06224'
$ws.Range("E50").Value = '  +0.01%  '

# Row 51
Set-TextValue "D51" '65.72'
$ws.Range("E51").Value = '  +2.29%  '
